$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 140.01849556838636
$ws.Range("C2").Value = 182.55403063177212
$ws.Range("D2").Value = 136.13812285768728
$ws.Range("E2").Value = 180.54888624821382

$ws.Range("B3").Value = 130.37329613799668
$ws.Range("C3").Value = 182.16845230034184
$ws.Range("D3").Value = 132.86164789905024
$ws.Range("E3").Value = 177.53767444996441

$ws.Range("B1:E3").Select()
